$xlPasteFormats = -4122

function SetIndex($ws, $r, $val) {
    $ws.Cells.Item($r, 1).Value = $val
}

function SetText($ws, $r, $c, $val) {
    # Leading apostrophe forces Excel to treat a numeric-looking string as
    # text (preserves leading zeros / trailing decimal zeros). The stray
    # "quotePrefix" style this creates is cleaned up later in bulk.
    $ws.Cells.Item($r, $c).Value = "'" + $val
}

function SetRank($ws, $r, $val) {
    $ws.Cells.Item($r, 8).Value = $val
}

$wb = $excel.ActiveWorkbook
$totalSheet = $wb.Worksheets.Item(1)

# ============================================================
# Part 1: insert the new 2022-Q3 row into the "总计" (total) sheet
# ============================================================

# Push existing rows 2..8 down to 3..9, duplicating row 2's formatting.
$totalSheet.Rows.Item(2).Insert()

# New row 2 holds the 2022-Q3 summary figures.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q3"
$totalSheet.Cells.Item(2, 3).Value = 35
$totalSheet.Cells.Item(2, 4).Value = 28.29

# Insert() leaves a stray style on B2:D2 and doesn't carry A2's original
# bold/bordered index style - fix both from neighbouring "clean" cells.
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B2:D2").PasteSpecial($xlPasteFormats)
$totalSheet.Range("A3").Copy()
$totalSheet.Range("A2").PasteSpecial($xlPasteFormats)

# Re-number the running index in column A for the rows that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $totalSheet.Cells.Item($r, 1).Value = $r - 2
}

# ============================================================
# Part 2: add the new "2022-Q3" worksheet right after "总计"
# ============================================================

$newSheet = $wb.Worksheets.Add($null, $totalSheet)
$newSheet.Name = "2022-Q3"

# Header row.
$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

SetIndex $newSheet 2 0
SetText  $newSheet 2 2 "501054"
SetText  $newSheet 2 3 "东方红睿泽三年定期开放灵活配置混合A"
SetText  $newSheet 2 4 "100.44"
SetText  $newSheet 2 5 "95.37"
SetText  $newSheet 2 6 "7.40"
SetText  $newSheet 2 7 "7.4326"
SetRank  $newSheet 2 1
SetIndex $newSheet 3 1
SetText  $newSheet 3 2 "009576"
SetText  $newSheet 3 3 "东方红智远三年持有期混合"
SetText  $newSheet 3 4 "62.81"
SetText  $newSheet 3 5 "92.83"
SetText  $newSheet 3 6 "7.09"
SetText  $newSheet 3 7 "4.4532"
SetRank  $newSheet 3 1
SetIndex $newSheet 4 2
SetText  $newSheet 4 2 "169104"
SetText  $newSheet 4 3 "东方红睿满沪港深灵活配置混合（LOF）"
SetText  $newSheet 4 4 "42.07"
SetText  $newSheet 4 5 "91.66"
SetText  $newSheet 4 6 "8.02"
SetText  $newSheet 4 7 "3.3740"
SetRank  $newSheet 4 1
SetIndex $newSheet 5 3
SetText  $newSheet 5 2 "000619"
SetText  $newSheet 5 3 "东方红产业升级灵活配置混合"
SetText  $newSheet 5 4 "40.04"
SetText  $newSheet 5 5 "93.62"
SetText  $newSheet 5 6 "7.08"
SetText  $newSheet 5 7 "2.8348"
SetRank  $newSheet 5 4
SetIndex $newSheet 6 4
SetText  $newSheet 6 2 "450009"
SetText  $newSheet 6 3 "国富中小盘股票"
SetText  $newSheet 6 4 "35.14"
SetText  $newSheet 6 5 "85.39"
SetText  $newSheet 6 6 "5.62"
SetText  $newSheet 6 7 "1.9749"
SetRank  $newSheet 6 1
SetIndex $newSheet 7 5
SetText  $newSheet 7 2 "000118"
SetText  $newSheet 7 3 "广发聚鑫债券A"
SetText  $newSheet 7 4 "145.84"
SetText  $newSheet 7 5 "20.01"
SetText  $newSheet 7 6 "0.98"
SetText  $newSheet 7 7 "1.4292"
SetRank  $newSheet 7 10
SetIndex $newSheet 8 6
SetText  $newSheet 8 2 "450002"
SetText  $newSheet 8 3 "国富弹性市值混合"
SetText  $newSheet 8 4 "30.14"
SetText  $newSheet 8 5 "85.83"
SetText  $newSheet 8 6 "4.70"
SetText  $newSheet 8 7 "1.4166"
SetRank  $newSheet 8 4
SetIndex $newSheet 9 7
SetText  $newSheet 9 2 "501066"
SetText  $newSheet 9 3 "东方红恒元五年定期开放灵活配置混合"
SetText  $newSheet 9 4 "19.71"
SetText  $newSheet 9 5 "96.96"
SetText  $newSheet 9 6 "5.15"
SetText  $newSheet 9 7 "1.0151"
SetRank  $newSheet 9 5
SetIndex $newSheet 10 8
SetText  $newSheet 10 2 "001409"
SetText  $newSheet 10 3 "工银互联网加股票"
SetText  $newSheet 10 4 "34.32"
SetText  $newSheet 10 5 "80.51"
SetText  $newSheet 10 6 "2.61"
SetText  $newSheet 10 7 "0.8958"
SetRank  $newSheet 10 10
SetIndex $newSheet 11 9
SetText  $newSheet 11 2 "011152"
SetText  $newSheet 11 3 "富兰克林国海兴海回报混合"
SetText  $newSheet 11 4 "14.87"
SetText  $newSheet 11 5 "84.99"
SetText  $newSheet 11 6 "3.92"
SetText  $newSheet 11 7 "0.5829"
SetRank  $newSheet 11 8
SetIndex $newSheet 12 10
SetText  $newSheet 12 2 "009121"
SetText  $newSheet 12 3 "广发招享混合A"
SetText  $newSheet 12 4 "58.31"
SetText  $newSheet 12 5 "23.12"
SetText  $newSheet 12 6 "0.94"
SetText  $newSheet 12 7 "0.5481"
SetRank  $newSheet 12 10
SetIndex $newSheet 13 11
SetText  $newSheet 13 2 "011468"
SetText  $newSheet 13 3 "国富竞争优势三年持有期混合A"
SetText  $newSheet 13 4 "11.65"
SetText  $newSheet 13 5 "82.69"
SetText  $newSheet 13 6 "3.41"
SetText  $newSheet 13 7 "0.3973"
SetRank  $newSheet 13 8
SetIndex $newSheet 14 12
SetText  $newSheet 14 2 "011069"
SetText  $newSheet 14 3 "工银成长精选混合A"
SetText  $newSheet 14 4 "12.40"
SetText  $newSheet 14 5 "60.10"
SetText  $newSheet 14 6 "2.17"
SetText  $newSheet 14 7 "0.2691"
SetRank  $newSheet 14 9
SetIndex $newSheet 15 13
SetText  $newSheet 15 2 "487021"
SetText  $newSheet 15 3 "工银优质精选混合"
SetText  $newSheet 15 4 "8.54"
SetText  $newSheet 15 5 "76.59"
SetText  $newSheet 15 6 "3.15"
SetText  $newSheet 15 7 "0.2690"
SetRank  $newSheet 15 4
SetIndex $newSheet 16 14
SetText  $newSheet 16 2 "013880"
SetText  $newSheet 16 3 "广发招享混合C"
SetText  $newSheet 16 4 "25.35"
SetText  $newSheet 16 5 "23.12"
SetText  $newSheet 16 6 "0.94"
SetText  $newSheet 16 7 "0.2383"
SetRank  $newSheet 16 10
SetIndex $newSheet 17 15
SetText  $newSheet 17 2 "009564"
SetText  $newSheet 17 3 "汇安消费龙头混合A"
SetText  $newSheet 17 4 "7.60"
SetText  $newSheet 17 5 "93.84"
SetText  $newSheet 17 6 "2.94"
SetText  $newSheet 17 7 "0.2234"
SetRank  $newSheet 17 10
SetIndex $newSheet 18 16
SetText  $newSheet 18 2 "000119"
SetText  $newSheet 18 3 "广发聚鑫债券C"
SetText  $newSheet 18 4 "16.95"
SetText  $newSheet 18 5 "20.01"
SetText  $newSheet 18 6 "0.98"
SetText  $newSheet 18 7 "0.1661"
SetRank  $newSheet 18 10
SetIndex $newSheet 19 17
SetText  $newSheet 19 2 "001496"
SetText  $newSheet 19 3 "工银聚焦30股票"
SetText  $newSheet 19 4 "3.56"
SetText  $newSheet 19 5 "84.46"
SetText  $newSheet 19 6 "4.49"
SetText  $newSheet 19 7 "0.1598"
SetRank  $newSheet 19 5
SetIndex $newSheet 20 18
SetText  $newSheet 20 2 "450010"
SetText  $newSheet 20 3 "国富策略回报混合"
SetText  $newSheet 20 4 "7.85"
SetText  $newSheet 20 5 "79.43"
SetText  $newSheet 20 6 "1.99"
SetText  $newSheet 20 7 "0.1562"
SetRank  $newSheet 20 7
SetIndex $newSheet 21 19
SetText  $newSheet 21 2 "000763"
SetText  $newSheet 21 3 "工银新财富灵活配置混合"
SetText  $newSheet 21 4 "2.76"
SetText  $newSheet 21 5 "92.61"
SetText  $newSheet 21 6 "5.37"
SetText  $newSheet 21 7 "0.1482"
SetRank  $newSheet 21 2
SetIndex $newSheet 22 20
SetText  $newSheet 22 2 "012880"
SetText  $newSheet 22 3 "国泰景气优选混合A"
SetText  $newSheet 22 4 "3.32"
SetText  $newSheet 22 5 "87.98"
SetText  $newSheet 22 6 "2.06"
SetText  $newSheet 22 7 "0.0684"
SetRank  $newSheet 22 8
SetIndex $newSheet 23 21
SetText  $newSheet 23 2 "020023"
SetText  $newSheet 23 3 "国泰事件驱动策略混合A"
SetText  $newSheet 23 4 "2.59"
SetText  $newSheet 23 5 "92.18"
SetText  $newSheet 23 6 "2.07"
SetText  $newSheet 23 7 "0.0536"
SetRank  $newSheet 23 9
SetIndex $newSheet 24 22
SetText  $newSheet 24 2 "011070"
SetText  $newSheet 24 3 "工银成长精选混合C"
SetText  $newSheet 24 4 "1.68"
SetText  $newSheet 24 5 "60.10"
SetText  $newSheet 24 6 "2.17"
SetText  $newSheet 24 7 "0.0365"
SetRank  $newSheet 24 9
SetIndex $newSheet 25 23
SetText  $newSheet 25 2 "001276"
SetText  $newSheet 25 3 "建信新经济灵活配置混合"
SetText  $newSheet 25 4 "1.13"
SetText  $newSheet 25 5 "70.13"
SetText  $newSheet 25 6 "3.22"
SetText  $newSheet 25 7 "0.0364"
SetRank  $newSheet 25 7
SetIndex $newSheet 26 24
SetText  $newSheet 26 2 "009954"
SetText  $newSheet 26 3 "北信瑞丰优选成长股票"
SetText  $newSheet 26 4 "0.55"
SetText  $newSheet 26 5 "93.35"
SetText  $newSheet 26 6 "4.43"
SetText  $newSheet 26 7 "0.0244"
SetRank  $newSheet 26 9
SetIndex $newSheet 27 25
SetText  $newSheet 27 2 "011032"
SetText  $newSheet 27 3 "东方红睿泽三年定期开放灵活配置混合C"
SetText  $newSheet 27 4 "0.32"
SetText  $newSheet 27 5 "95.37"
SetText  $newSheet 27 6 "7.40"
SetText  $newSheet 27 7 "0.0237"
SetRank  $newSheet 27 1
SetIndex $newSheet 28 26
SetText  $newSheet 28 2 "011469"
SetText  $newSheet 28 3 "国富竞争优势三年持有期混合C"
SetText  $newSheet 28 4 "0.67"
SetText  $newSheet 28 5 "82.69"
SetText  $newSheet 28 6 "3.41"
SetText  $newSheet 28 7 "0.0228"
SetRank  $newSheet 28 8
SetIndex $newSheet 29 27
SetText  $newSheet 29 2 "009565"
SetText  $newSheet 29 3 "汇安消费龙头混合C"
SetText  $newSheet 29 4 "0.49"
SetText  $newSheet 29 5 "93.84"
SetText  $newSheet 29 6 "2.94"
SetText  $newSheet 29 7 "0.0144"
SetRank  $newSheet 29 10
SetIndex $newSheet 30 28
SetText  $newSheet 30 2 "001829"
SetText  $newSheet 30 3 "北信瑞丰中国智造主题灵活配置混合"
SetText  $newSheet 30 4 "0.28"
SetText  $newSheet 30 5 "93.27"
SetText  $newSheet 30 6 "5.00"
SetText  $newSheet 30 7 "0.0140"
SetRank  $newSheet 30 7
SetIndex $newSheet 31 29
SetText  $newSheet 31 2 "012881"
SetText  $newSheet 31 3 "国泰景气优选混合C"
SetText  $newSheet 31 4 "0.22"
SetText  $newSheet 31 5 "87.98"
SetText  $newSheet 31 6 "2.06"
SetText  $newSheet 31 7 "0.0045"
SetRank  $newSheet 31 8
SetIndex $newSheet 32 30
SetText  $newSheet 32 2 "000761"
SetText  $newSheet 32 3 "国富健康优质生活股票"
SetText  $newSheet 32 4 "0.14"
SetText  $newSheet 32 5 "79.95"
SetText  $newSheet 32 6 "3.13"
SetText  $newSheet 32 7 "0.0044"
SetRank  $newSheet 32 8
SetIndex $newSheet 33 31
SetText  $newSheet 33 2 "003685"
SetText  $newSheet 33 3 "汇安丰融灵活配置混合C"
SetText  $newSheet 33 4 "0.09"
SetText  $newSheet 33 5 "93.48"
SetText  $newSheet 33 6 "3.45"
SetText  $newSheet 33 7 "0.0031"
SetRank  $newSheet 33 10
SetIndex $newSheet 34 32
SetText  $newSheet 34 2 "161718"
SetText  $newSheet 34 3 "招商沪深300高贝塔指数"
SetText  $newSheet 34 4 "0.15"
SetText  $newSheet 34 5 "94.52"
SetText  $newSheet 34 6 "1.50"
SetText  $newSheet 34 7 "0.0022"
SetRank  $newSheet 34 2
SetIndex $newSheet 35 33
SetText  $newSheet 35 2 "003684"
SetText  $newSheet 35 3 "汇安丰融灵活配置混合A"
SetText  $newSheet 35 4 "0.02"
SetText  $newSheet 35 5 "93.48"
SetText  $newSheet 35 6 "3.45"
SetText  $newSheet 35 7 "0.0007"
SetRank  $newSheet 35 10
SetIndex $newSheet 36 34
SetText  $newSheet 36 2 "015592"
SetText  $newSheet 36 3 "国泰事件驱动策略混合C"
SetText  $newSheet 36 4 "0.02"
SetText  $newSheet 36 5 "92.18"
SetText  $newSheet 36 6 "2.07"
SetText  $newSheet 36 7 "0.0004"
SetRank  $newSheet 36 9

# ============================================================
# Part 3: fix up styles on the new sheet in bulk
# ============================================================

# Header row B1:H1 should look like the bold/bordered header used on every
# other quarter sheet - copy it straight from the "总计" sheet's header.
$totalSheet.Range("B1").Copy()
$newSheet.Range("B1:H1").PasteSpecial($xlPasteFormats)

# Column A (the running index) uses the same bold/bordered style too.
$totalSheet.Range("A3").Copy()
$newSheet.Range("A2:A36").PasteSpecial($xlPasteFormats)

# Columns B:G got a stray quotePrefix style from the leading-apostrophe
# trick used by SetText - reset them to the plain default style used by
# every other data cell (e.g. the rank column H, which was never touched).
$newSheet.Range("H2").Copy()
$newSheet.Range("B2:G36").PasteSpecial($xlPasteFormats)

Write-Output "2022-Q3 sheet created and 总计 sheet updated"
